# A new daily price record (week of 2023-03-31) was inserted into the
# "Ciruela" log at row 69 of the active sheet, pushing every existing
# record from row 69 onward down by one row (old row 112 becomes row 113).
#
# Excel's native row-insert semantics give us exactly that shift (formats,
# dimension, everything below moves down) so we insert a blank row at 69
# and then populate it with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 69; rows 69-112 shift down to 70-113.
$ws.Rows.Item(69).Insert()

$newRow = 69

$ws.Cells.Item($newRow, 1).Value  = 7
$ws.Cells.Item($newRow, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($newRow, 3).Value  = "Ñuble"

$ws.Cells.Item($newRow, 4).Value  = 45016
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow + 1, 4).NumberFormat

$ws.Cells.Item($newRow, 5).Value  = 16
$ws.Cells.Item($newRow, 6).Value  = "Fruta"
$ws.Cells.Item($newRow, 7).Value  = 100103
$ws.Cells.Item($newRow, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item($newRow, 9).Value  = 100103002
$ws.Cells.Item($newRow, 10).Value = "Ciruela"
$ws.Cells.Item($newRow, 11).Value = "Angeleno"
$ws.Cells.Item($newRow, 12).Value = "Primera"
$ws.Cells.Item($newRow, 13).Value = 40
$ws.Cells.Item($newRow, 14).Value = 12000
$ws.Cells.Item($newRow, 15).Value = 12000
$ws.Cells.Item($newRow, 16).Value = 12000
$ws.Cells.Item($newRow, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item($newRow, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($newRow, 19).Value = 667
$ws.Cells.Item($newRow, 20).Value = 18
